# Sheet 5 ("保險" / insurance) had two bugs in the exported data:
#   1. Row 1 (the header row) was accidentally populated with the first
#      data row's values instead of real header labels.
#   2. Only 4 columns (A-D) were exported; the common trailing columns
#      that every other sheet has (property_category, category, date,
#      legislator_name, legislator_id, source_file, index) were missing.
# This script fixes the header row and appends the missing columns for
# both existing data rows, matching the layout used by the other sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# ---- Header row (row 1): replace the stray data values with real
# ---- header labels, and add headers for the new trailing columns.
$ws.Cells.Item(1,2).Value = "company"
$ws.Cells.Item(1,3).Value = "name"
$ws.Cells.Item(1,4).Value = "owner"
$ws.Cells.Item(1,5).Value = "property_category"
$ws.Cells.Item(1,6).Value = "category"
$ws.Cells.Item(1,7).Value = "date"
$ws.Cells.Item(1,8).Value = "legislator_name"
$ws.Cells.Item(1,9).Value = "legislator_id"
$ws.Cells.Item(1,10).Value = "source_file"
$ws.Cells.Item(1,11).Value = "index"

# Give the newly added header cells (E1:K1) the same bold/bordered/
# centered formatting already used by B1:D1.
$ws.Cells.Item(1,4).Copy()
$ws.Range($ws.Cells.Item(1,5), $ws.Cells.Item(1,11)).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Row 2 data (index 80): company/name were shifted over; fill in
# ---- the newly added columns.
$ws.Cells.Item(2,2).Value = "新光人壽"
$ws.Cells.Item(2,3).Value = "新光人壽全心全意終身還本保險"
$ws.Cells.Item(2,5).Value = "insurance"
$ws.Cells.Item(2,6).Value = "normal"
# Assign the date through a text formula so Excel doesn't reinterpret the
# ISO-looking string as a date serial, then freeze it down to a plain value.
$ws.Cells.Item(2,7).Formula = '="2012-02-29"'
$ws.Cells.Item(2,7).Copy()
$ws.Cells.Item(2,7).PasteSpecial(-4163)
$ws.Cells.Item(2,8).Value = "林岱樺"
$ws.Cells.Item(2,9).Value = 904
$ws.Cells.Item(2,10).Value = "tmp3bff1"
$ws.Cells.Item(2,11).Value = 80

# ---- Row 3 data (index 81): name shifted over; fill in the new columns.
$ws.Cells.Item(3,3).Value = "世紀領航萬能終身壽險計劃A"
$ws.Cells.Item(3,5).Value = "insurance"
$ws.Cells.Item(3,6).Value = "normal"
$ws.Cells.Item(3,7).Formula = '="2012-02-29"'
$ws.Cells.Item(3,7).Copy()
$ws.Cells.Item(3,7).PasteSpecial(-4163)
$ws.Cells.Item(3,8).Value = "林岱樺"
$ws.Cells.Item(3,9).Value = 904
$ws.Cells.Item(3,10).Value = "tmp3bff1"
$ws.Cells.Item(3,11).Value = 81

# Give the new data cells in both rows (E:K) the same plain formatting
# already used by the existing data cells (A:D) in those rows.
$ws.Cells.Item(2,4).Copy()
$ws.Range($ws.Cells.Item(2,5), $ws.Cells.Item(2,11)).PasteSpecial(-4122)
$ws.Cells.Item(3,4).Copy()
$ws.Range($ws.Cells.Item(3,5), $ws.Cells.Item(3,11)).PasteSpecial(-4122)
$excel.CutCopyMode = $false
